$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Saldo" column (B2:B11) now shows amounts in the same Euro accounting format
# already used by the other money columns (G:I).
$ws.Range("B2:B11").NumberFormat = "_-* #,##0.00\ [$€-410]_-;\-* #,##0.00\ [$€-410]_-;_-* ""-""??\ [$€-410]_-;_-@_-"

# "Acconti richiesti" advance payment for apartment 1 updated from -1800 to -2000.
$ws.Range("G12").Value = -2000

# Touch the otherwise-blank cells in the "Acconti richiesti" / "Conguaglio spese"
# total rows so they carry the row's own formatting (matches the rest of the row).
$ws.Range("B12:F12").Font.Bold = $ws.Range("B12:F12").Font.Bold
$ws.Range("B13:F13").Font.Bold = $ws.Range("B13:F13").Font.Bold
